$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Object_Mapping")

# Rename Solar_Plant_Kasso -> Solar_Plant
$ws.Range("A2").Value = "Solar_Plant"

# Insert a new row before row 8 (power_line_Wholesale_Kasso) for the Wind Farm
$ws.Rows.Item(8).Insert()

# Populate the new row 8 with Wind Farm data
$ws.Range("A8").Value = "Wind_Farm"
$ws.Range("B8").Value = "Wind_farm"
